# Auto-generated Excel COM-interop script
# Applies numeric corrections to the crafting-leve profit tables
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.
# (scheduled data-refresh run: currentAveragePrice* / LevePrice* / LeveProfit* columns)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2755.3
$ws.Range("I100").Value = 499
$ws.Range("J100").Value = 3006
$ws.Range("K100").Value = 499
$ws.Range("L100").Value = 3006
$ws.Range("M100").Value = 42
$ws.Range("N100").Value = -4088
$ws.Range("H115").Value = 4147.5
$ws.Range("I115").Value = 5245.8335
$ws.Range("J115").Value = 2500
$ws.Range("K115").Value = 15737.5005
$ws.Range("L115").Value = 7500
$ws.Range("M115").Value = -14170.5005
$ws.Range("N115").Value = -10634
$ws.Range("H127").Value = 1239.5
$ws.Range("I127").Value = 887.4
$ws.Range("J127").Value = 3000
$ws.Range("K127").Value = 2662.2
$ws.Range("L127").Value = 9000
$ws.Range("M127").Value = 2297.8
$ws.Range("N127").Value = -18920
$ws.Range("H138").Value = 3519.162
$ws.Range("I138").Value = 1292.2632
$ws.Range("J138").Value = 4288.4546
$ws.Range("K138").Value = 3876.7896
$ws.Range("L138").Value = 12865.3638
$ws.Range("M138").Value = 1263.2104
$ws.Range("N138").Value = -23145.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7780.122
$ws.Range("I32").Value = 5780.406
$ws.Range("J32").Value = 18394
$ws.Range("K32").Value = 5780.406
$ws.Range("L32").Value = 18394
$ws.Range("M32").Value = -5493.406
$ws.Range("N32").Value = -18968
$ws.Range("H74").Value = 5325.6216
$ws.Range("I74").Value = 2771.1333
$ws.Range("J74").Value = 16273.429
$ws.Range("K74").Value = 2771.1333
$ws.Range("L74").Value = 16273.429
$ws.Range("M74").Value = -1897.1333
$ws.Range("N74").Value = -18021.429
$ws.Range("H77").Value = 5325.6216
$ws.Range("I77").Value = 2771.1333
$ws.Range("J77").Value = 16273.429
$ws.Range("K77").Value = 13855.6665
$ws.Range("L77").Value = 81367.145
$ws.Range("M77").Value = -9487.666499999999
$ws.Range("N77").Value = -90103.145
$ws.Range("H102").Value = 2797.6924
$ws.Range("I102").Value = 2595.7144
$ws.Range("J102").Value = 3033.3333
$ws.Range("K102").Value = 2595.7144
$ws.Range("L102").Value = 3033.3333
$ws.Range("M102").Value = -973.7143999999998
$ws.Range("N102").Value = -6277.3333
$ws.Range("H132").Value = 5719.795
$ws.Range("I132").Value = 1356.3334
$ws.Range("K132").Value = 4069.0002
$ws.Range("M132").Value = -1539.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 5670.9
$ws.Range("I36").Value = 1878.7778
$ws.Range("J36").Value = 39800
$ws.Range("K36").Value = 1878.7778
$ws.Range("L36").Value = 39800
$ws.Range("M36").Value = -1344.7778
$ws.Range("N36").Value = -40868
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 5632.839
$ws.Range("I105").Value = 5379.9
$ws.Range("J105").Value = 6092.727
$ws.Range("K105").Value = 5379.9
$ws.Range("L105").Value = 6092.727
$ws.Range("M105").Value = -3632.9
$ws.Range("N105").Value = -9586.726999999999
$ws.Range("H111").Value = 57766.668
$ws.Range("J111").Value = 57766.668
$ws.Range("L111").Value = 57766.668
$ws.Range("N111").Value = -65946.66800000001
$ws.Range("H133").Value = 67136
$ws.Range("J133").Value = 67136
$ws.Range("L133").Value = 67136
$ws.Range("N133").Value = -77256

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7933.8335
$ws.Range("I4").Value = 1500
$ws.Range("J4").Value = 11150.75
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 11150.75
$ws.Range("M4").Value = -1388
$ws.Range("N4").Value = -11374.75
$ws.Range("H103").Value = 13006.857
$ws.Range("I103").Value = 9341.333000000001
$ws.Range("K103").Value = 9341.333000000001
$ws.Range("M103").Value = -8169.333000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8389.615
$ws.Range("I68").Value = 507.2857
$ws.Range("J68").Value = 17585.666
$ws.Range("K68").Value = 1521.8571
$ws.Range("L68").Value = 52756.99800000001
$ws.Range("M68").Value = -710.8571000000002
$ws.Range("N68").Value = -54378.99800000001
$ws.Range("H71").Value = 8389.615
$ws.Range("I71").Value = 507.2857
$ws.Range("J71").Value = 17585.666
$ws.Range("K71").Value = 4565.571300000001
$ws.Range("L71").Value = 158270.994
$ws.Range("M71").Value = -509.5713000000005
$ws.Range("N71").Value = -166382.994
$ws.Range("H75").Value = 3146.3635
$ws.Range("I75").Value = 750
$ws.Range("J75").Value = 3678.889
$ws.Range("K75").Value = 2250
$ws.Range("L75").Value = 11036.667
$ws.Range("M75").Value = -1252
$ws.Range("N75").Value = -13032.667
$ws.Range("H78").Value = 3146.3635
$ws.Range("I78").Value = 750
$ws.Range("J78").Value = 3678.889
$ws.Range("K78").Value = 6750
$ws.Range("L78").Value = 33110.001
$ws.Range("M78").Value = -1758
$ws.Range("N78").Value = -43094.001
$ws.Range("H87").Value = 8387.684999999999
$ws.Range("I87").Value = 4602.8
$ws.Range("J87").Value = 9739.429
$ws.Range("K87").Value = 13808.4
$ws.Range("L87").Value = 29218.287
$ws.Range("M87").Value = -12560.4
$ws.Range("N87").Value = -31714.287
$ws.Range("H90").Value = 8387.684999999999
$ws.Range("I90").Value = 4602.8
$ws.Range("J90").Value = 9739.429
$ws.Range("K90").Value = 41425.2
$ws.Range("L90").Value = 87654.861
$ws.Range("M90").Value = -35185.2
$ws.Range("N90").Value = -100134.861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 223.42857
$ws.Range("I2").Value = 94
$ws.Range("J2").Value = 396
$ws.Range("K2").Value = 94
$ws.Range("L2").Value = 396
$ws.Range("M2").Value = 19
$ws.Range("N2").Value = -622
$ws.Range("H46").Value = 15176.333
$ws.Range("I46").Value = 10041
$ws.Range("J46").Value = 15818.25
$ws.Range("K46").Value = 10041
$ws.Range("L46").Value = 15818.25
$ws.Range("M46").Value = -9885
$ws.Range("N46").Value = -16130.25
$ws.Range("H132").Value = 7773.3335
$ws.Range("I132").Value = 11495
$ws.Range("J132").Value = 3121.25
$ws.Range("K132").Value = 34485
$ws.Range("L132").Value = 9363.75
$ws.Range("M132").Value = -31955
$ws.Range("N132").Value = -14423.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3744.7407
$ws.Range("I40").Value = 3582.182
$ws.Range("J40").Value = 4460
$ws.Range("K40").Value = 3582.182
$ws.Range("L40").Value = 4460
$ws.Range("M40").Value = -3446.182
$ws.Range("N40").Value = -4732
$ws.Range("H68").Value = 2237.5
$ws.Range("I68").Value = 1820
$ws.Range("K68").Value = 1820
$ws.Range("M68").Value = -1071
$ws.Range("H71").Value = 2237.5
$ws.Range("I71").Value = 1820
$ws.Range("K71").Value = 9100
$ws.Range("M71").Value = -5356
$ws.Range("H136").Value = 8581.143
$ws.Range("I136").Value = 9400.666999999999
$ws.Range("J136").Value = 8253.333000000001
$ws.Range("K136").Value = 28202.001
$ws.Range("L136").Value = 24759.999
$ws.Range("M136").Value = -25652.001
$ws.Range("N136").Value = -29859.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1552.2941
$ws.Range("I132").Value = 1399.0769
$ws.Range("J132").Value = 2050.25
$ws.Range("K132").Value = 4197.2307
$ws.Range("L132").Value = 6150.75
$ws.Range("M132").Value = -1667.2307
$ws.Range("N132").Value = -11210.75
